$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header texts (keep A1 "Nombre" as-is)
$ws.Range("B1").Value = "Año (aaaa)"
$ws.Range("C1").Value = "Fecha de inicio de inscripciones (dd/mm/aaaa)"
$ws.Range("D1").Value = "Fecha de fin de inscripciones (dd/mm/aaaa)"
$ws.Range("E1").Value = "Fecha fin de oferta (dd/mm/aaaa)"

# Move/update the active selection to B1
$ws.Range("B1").Select()
